$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Relacionamento")
$ws.Range("A1").Value = "test"
